$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with GitHub Actions scraped data.
# Force column D (Price) to text format before writing so that
# numeric-looking values (e.g. "243.85") are not auto-converted
# to floating point numbers by Excel, keeping them as text like the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "44.256.47"
$ws.Range("E2").Value = "  +1.67%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "2.364.46"
$ws.Range("E3").Value = "  -0.88%  "
# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.03%  "
# Row 5: XRP
$ws.Range("E5").Value = "  +5.14%  "
# Row 6: BNB
$ws.Range("D6").Value = "243.85"
$ws.Range("E6").Value = "  +3.12%  "
# Row 7: Solana
$ws.Range("D7").Value = "74.14"
$ws.Range("E7").Value = "  +2.56%  "
# Row 8: USDC
$ws.Range("E8").Value = "  -0.04%  "
# Row 9: Cardano
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +28.25%  "
# Row 10: Dogecoin
$ws.Range("E10").Value = "  +5.95%  "
# Row 11: Avalanche
$ws.Range("E11").Value = "  +16.75%  "
# Row 12: Polkadot
$ws.Range("D12").Value = "7.55"
$ws.Range("E12").Value = "  +20.00%  "
# Row 13: TRON
$ws.Range("E13").Value = "  +2.07%  "
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.718.03"
$ws.Range("E14").Value = "  -0.79%  "
# Row 15: Chainlink
$ws.Range("D15").Value = "17.01"
$ws.Range("E15").Value = "  +6.70%  "
# Row 16: Polygon
$ws.Range("D16").Value = "0.918"
$ws.Range("E16").Value = "  +6.84%  "
# Row 17: WrappedEther
$ws.Range("D17").Value = "2.357.16"
$ws.Range("E17").Value = "  -1.25%  "
# Row 18: WrappedBTC
$ws.Range("D18").Value = "44.331.42"
$ws.Range("E18").Value = "  +1.85%  "
# Row 19: ShibaInu
$ws.Range("E19").Value = "  +4.75%  "
# Row 20: Uniswap
$ws.Range("D20").Value = "6.76"
$ws.Range("E20").Value = "  +5.57%  "
# Row 21: Litecoin
$ws.Range("D21").Value = "78.59"
$ws.Range("E21").Value = "  +5.21%  "
# Row 22: BitcoinCash
$ws.Range("D22").Value = "257.08"
$ws.Range("E22").Value = "  +2.09%  "
# Row 23: Dai
$ws.Range("E23").Value = "  +0.10%  "
# Row 24: PancakeSwap
$ws.Range("E24").Value = "  +3.57%  "
# Row 25: WEMIXToken
$ws.Range("E25").Value = "  -2.74%  "
# Row 26: Cosmos
$ws.Range("D26").Value = "10.84"
$ws.Range("E26").Value = "  +7.91%  "
# Row 27: Toncoin
$ws.Range("E27").Value = "  +1.35%  "
# Row 28: EthereumClassic
$ws.Range("D28").Value = "22.73"
$ws.Range("E28").Value = "  -1.08%  "
# Row 29: ImmutableX
$ws.Range("D29").Value = "1.64"
$ws.Range("E29").Value = "  +6.50%  "
# Row 30: Monero
$ws.Range("D30").Value = "175.29"
$ws.Range("E30").Value = "  +0.43%  "
# Row 31: Kaspa
$ws.Range("E31").Value = "  +1.64%  "
# Row 32: Stellar
$ws.Range("E32").Value = "  +5.88%  "
# Row 33: Filecoin
$ws.Range("D33").Value = "5.42"
$ws.Range("E33").Value = "  +8.07%  "
# Row 34: Hedera
$ws.Range("D34").Value = "0.0763"
$ws.Range("E34").Value = "  +9.89%  "
# Row 35: InternetComputer(DFINITY)
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  +6.30%  "
# Row 36: RenderToken
$ws.Range("E36").Value = "  +4.90%  "
# Row 37: LidoDAOToken
$ws.Range("D37").Value = "2.48"
$ws.Range("E37").Value = "  +0.42%  "
# Row 38: THORChain
$ws.Range("D38").Value = "6.55"
$ws.Range("E38").Value = "  -1.27%  "
# Row 39: VeChain
$ws.Range("D39").Value = "0.0275"
$ws.Range("E39").Value = "  +7.18%  "
# Row 40: InjectiveProtocol
$ws.Range("E40").Value = "  +0.04%  "
# Row 41: FraxShare
$ws.Range("E41").Value = "  +1.51%  "
# Row 42: BinanceUSD
$ws.Range("E42").Value = "  +0.00%  "
# Row 43: Algorand
$ws.Range("D43").Value = "0.198"
$ws.Range("E43").Value = "  +17.01%  "
# Row 44: Cronos
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.100"
$ws.Range("E44").Value = "  +5.20%  "
# Row 45: NEARProtocol
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +11.18%  "
# Row 46: TrustWalletToken
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.26"
$ws.Range("E46").Value = "  +2.74%  "
# Row 47: Aave
$ws.Range("D47").Value = "101.54"
$ws.Range("E47").Value = "  +1.43%  "
# Row 48: ARBITRUM
$ws.Range("E48").Value = "  -0.57%  "
# Row 49: FTXToken
$ws.Range("E49").Value = "  -2.54%  "
# Row 50: Maker
$ws.Range("D50").Value = "1.462.66"
$ws.Range("E50").Value = "  +0.70%  "
# Row 51: MultiversX
$ws.Range("D51").Value = "53.39"
$ws.Range("E51").Value = "  +4.75%  "

# Restore the default (General) cell formatting for column D now that
# the text values have been safely written, matching the workbook's original styling.
$ws.Range("D2:D51").ClearFormats()
